$wb = $excel.ActiveWorkbook

# New player "D.Helm" is logged as the next column after the last
# existing player (column R -> column S) on both the Rushing and
# Receiving sheets: a header cell matching the other header cells'
# formatting, and a data cell ("n", i.e. not yet simulated/no value)
# matching the other data cells' formatting.
foreach ($ws in $wb.Worksheets) {
    $ws.Range("S1").Value = "D.Helm"
    $ws.Range("R1").Copy()
    $ws.Range("S1").PasteSpecial(-4122)

    $ws.Range("S2").Value = "n"
    $ws.Range("R2").Copy()
    $ws.Range("S2").PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
